$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "-1" suffix from the city names used as Source/Destination labels.
# Using Range.Replace (rather than setting .Value on individual cells) edits
# the shared-string text in place, so every cell referencing that string
# (column B "Source" and column C "Destination") is updated together.
$ws.Cells.Replace("Tehran-1", "Tehran")
$ws.Cells.Replace("Mashhad-1", "Mashhad")
$ws.Cells.Replace("Esfahan-1", "Esfahan")
$ws.Cells.Replace("Shiraz-1", "Shiraz")
$ws.Cells.Replace("Ahwaz-1", "Ahwaz")
$ws.Cells.Replace("Tabriz-1", "Tabriz")
$ws.Cells.Replace("Babol-1", "Babol")
$ws.Cells.Replace("Hamedan-1", "Hamedan")

# Leave the cursor where the author finished editing (Destination cell of
# row 17, just after typing "Hamedan" into C17 and tabbing to D17).
$ws.Range("D17").Select()
